$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text so Excel does not silently convert them from strings to numbers
# (the source data stores every Price/Volume cell as text).

$ws.Range('D2').Value = '36.463.84'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '1.949.53'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.613'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.376'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('E10').Value = '  -7.33%  '
$ws.Range('E11').Value = '  -1.49%  '
$ws.Range('D12').Value = '2.234.63'
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.23'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').Value = '1.946.25'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').Value = '36.324.30'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('D20').Value = '0.0₃0846'
$ws.Range('E20').Value = '  -4.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '228.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('E24').Value = '  +2.32%  '
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.07'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E28').Value = '  +0.97%  '
$ws.Range('E29').Value = '  -0.79%  '
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.67'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0609'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.32'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.41'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +12.27%  '
$ws.Range('E37').Value = '  +3.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.76'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.21'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -15.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0970'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.40%  '
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('E42').Value = '  -1.22%  '
$ws.Range('E43').Value = '  -1.28%  '
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').Value = '1.360.68'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.82'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('D50').Value = '2.126.01'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.94%  '
